$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.631.91'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '3.073.73'
$ws.Range('E3').Value = '  -2.51%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''236.01'
$ws.Range('E5').Value = '  +9.36%  '
$ws.Range('D6').Value = '''618.43'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('E7').Value = '  -10.32%  '
$ws.Range('D8').Value = '''0.361'
$ws.Range('E8').Value = '  -1.46%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = '3.070.29'
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('E11').Value = '  -5.70%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '''0.0000249'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').Value = '''35.24'
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('D15').Value = '89.589.39'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').Value = '''5.36'
$ws.Range('E16').Value = '  -6.50%  '
$ws.Range('D17').Value = '3.658.32'
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('D18').Value = '3.100.68'
$ws.Range('E18').Value = '  -3.03%  '
$ws.Range('D19').Value = '''3.80'
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').Value = '''13.73'
$ws.Range('E21').Value = '  -6.20%  '
$ws.Range('D22').Value = '''433.06'
$ws.Range('E22').Value = '  -8.52%  '
$ws.Range('D23').Value = '''5.40'
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('D24').Value = '''8.76'
$ws.Range('E24').Value = '  -4.22%  '
$ws.Range('D25').Value = '''5.73'
$ws.Range('E25').Value = '  -0.46%  '
$ws.Range('D26').Value = '''86.85'
$ws.Range('E26').Value = '  -8.66%  '
$ws.Range('D27').Value = '''11.75'
$ws.Range('E27').Value = '  -5.13%  '
$ws.Range('D28').Value = '3.253.11'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('D29').Value = '''0.997'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '''9.05'
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('E32').Value = '  -4.12%  '
$ws.Range('E33').Value = '  -8.53%  '
$ws.Range('E34').Value = '  +4.55%  '
$ws.Range('D35').Value = '''25.56'
$ws.Range('E35').Value = '  -6.93%  '
$ws.Range('D36').Value = '''3.70'
$ws.Range('E36').Value = '  +2.72%  '
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('D38').Value = '''495.14'
$ws.Range('E38').Value = '  -4.48%  '
$ws.Range('E39').Value = '  -3.23%  '
$ws.Range('E40').Value = '  -3.57%  '
$ws.Range('D41').Value = '''0.0895'
$ws.Range('E41').Value = '  -2.29%  '
$ws.Range('D42').Value = '''3.60'
$ws.Range('E42').Value = '  +53.83%  '
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  -6.88%  '
$ws.Range('D46').Value = '''151.90'
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('E47').Value = '  -6.66%  '
$ws.Range('D48').Value = '''0.674'
$ws.Range('E48').Value = '  -9.76%  '
$ws.Range('D49').Value = '''44.35'
$ws.Range('E49').Value = '  -2.54%  '
$ws.Range('D50').Value = '''0.998'
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('E51').Value = '  -4.51%  '
